# Updated Master Gantt Chart
#
# Populates real task names / actual hours that were previously TBD
# placeholders on the Gantt sheet, records attendance + duration for the
# second ("Jan. 28") team meeting, and fills in the Systems-Analysis
# actual-hours tracker for every team member. The Management Summary sheet
# pulls from all three via formulas, so it recalculates on its own.
# Finally, re-points the active sheet/tab + per-sheet selections the way
# the author left them.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Gantt sheet: give the first task of 4 team members its real name, and
# fill in the hours budgeted / actual for the ones that have progressed.
# (Order matters here only insofar as it matches the order the strings
# were authored in; it has no functional effect.)
# ------------------------------------------------------------------
$gantt = $wb.Worksheets.Item("Gantt")

# Corbin Schueller - task 1
$gantt.Cells.Item(63, 1).Value = "Logo Design"
$gantt.Cells.Item(63, 2).Value = 1
$gantt.Cells.Item(63, 3).Value = 0.2

# Benjamin Hallman - task 1 (name only; hours still 0/0)
$gantt.Cells.Item(39, 1).Value = "Setup Master Gantt"

# Delaney Fitzgerald - task 1
$gantt.Cells.Item(51, 1).Value = "Coding Standards"
$gantt.Cells.Item(51, 2).Value = 1
$gantt.Cells.Item(51, 3).Value = 2

# Jacob Friedberg - task 1
$gantt.Cells.Item(3, 1).Value = "GitHub and Discord Setup"
$gantt.Cells.Item(3, 2).Value = 1
$gantt.Cells.Item(3, 3).Value = 0.5

# ------------------------------------------------------------------
# Meetings sheet: "Jan. 28" meeting (column D) ran half an hour, and
# Conrad, Benjamin, Delaney & Corbin attended it (Jacob and Cameron did
# not).
# ------------------------------------------------------------------
$meetings = $wb.Worksheets.Item("Meetings")

$meetings.Cells.Item(3, 4).Value = 0.5

$meetings.Cells.Item(6, 4).Value = "ü"
$meetings.Cells.Item(6, 4).HorizontalAlignment = -4108

$meetings.Cells.Item(7, 4).Value = "ü"
$meetings.Cells.Item(7, 4).HorizontalAlignment = -4108

$meetings.Cells.Item(8, 4).Value = "ü"
$meetings.Cells.Item(8, 4).HorizontalAlignment = -4108

$meetings.Cells.Item(9, 4).Value = "ü"
$meetings.Cells.Item(9, 4).HorizontalAlignment = -4108

# ------------------------------------------------------------------
# SA sheet: actual hours spent on Champion / RFP / SA Presentation for
# every team member (column C of each 3-row block).
# ------------------------------------------------------------------
$sa = $wb.Worksheets.Item("SA")

$saRows = 2, 6, 10, 14, 18, 22
foreach ($r in $saRows) {
    $sa.Cells.Item($r, 3).Value = 3
    $sa.Cells.Item($r + 1, 3).Value = 4
    $sa.Cells.Item($r + 2, 3).Value = 2
}

# ------------------------------------------------------------------
# Restore per-sheet selections, then land on Management Summary as the
# active tab (matching the saved workbook state).
# ------------------------------------------------------------------
$sa.Range("C17").Select()
$meetings.Range("D3").Select()
$gantt.Range("E17").Select()

$mgmt = $wb.Worksheets.Item("Management Summary")
$mgmt.Activate()
$mgmt.Range("J4").Select()
